$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename "Top 5" -> "Top 6" (this also updates any formulas that
#    reference the sheet, e.g. defined names / cell formulas -- the
#    chart's cached series formula is fixed up explicitly below).
# ------------------------------------------------------------------
$topSheet = $wb.Worksheets.Item("Top 5")
$topSheet.Name = "Top 6"

# Fix up the chart living on the renamed sheet so its series formula
# points at the new sheet name (cached chart XML isn't auto-rewritten
# by a plain rename in this engine).
$topChartObj = $topSheet.ChartObjects().Item(1)
$topSeries = $topChartObj.Chart.SeriesCollection().Item(1)
$topSeries.Formula = "=SERIES(,'Top 6'!`$B`$38:`$B`$43,'Top 6'!`$C`$38:`$C`$43,1)"

# Selection on the "Top 6" sheet moved to H22.
$topSheet.Range("H22").Select()

# ------------------------------------------------------------------
# 2. Add a new "# records" sheet after "feature importance" with the
#    isotope-combination record counts, and a bar chart of the data
#    (mirroring the existing "feature importance" sheet/chart).
# ------------------------------------------------------------------
$featSheet = $wb.Worksheets.Item("feature importance")
$recSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $featSheet)
$recSheet.Name = "# records"

$recSheet.Range("A1").Value = "isotope"
$recSheet.Range("A2").Value = "combinations"
$recSheet.Range("B2").Value = "# records"
$recSheet.Range("A2:B2").Font.Underline = 2

$recSheet.Range("A3").Value = "C Si N"
$recSheet.Range("B3").Value = 1301
$recSheet.Range("A4").Value = "C N"
$recSheet.Range("B4").Value = 2189
$recSheet.Range("A5").Value = "C Si"
$recSheet.Range("B5").Value = 14423
$recSheet.Range("B5").NumberFormat = "#,##0"

$recSheet.Columns.Item(1).ColumnWidth = 12.33203125

$recChartObj = $recSheet.ChartObjects().Add(60, 10, 320, 216)
$recChart = $recChartObj.Chart
$recChart.ChartType = 57
$recChart.SetSourceData($recSheet.Range("A2:B5"))
$recSeries = $recChart.SeriesCollection().Item(1)
$recSeries.HasDataLabels = $true

# This new chart becomes the most-recently-inserted "quick chart", so
# the workbook's hidden _xlchart.v2.* bookkeeping names now point at
# its source data instead of the "feature importance" chart's.
$wb.Names.Item("_xlchart.v2.0").RefersTo = "='# records'!`$A`$3:`$A`$5"
$wb.Names.Item("_xlchart.v2.1").RefersTo = "='# records'!`$B`$2"
$wb.Names.Item("_xlchart.v2.2").RefersTo = "='# records'!`$B`$3:`$B`$5"

# The new sheet becomes the active tab, with B1 selected.
$recSheet.Activate()
$recSheet.Range("B1").Select()
